# Ind1.xlsx - "added details on data structure"
#
# Three data cells on Sheet1 get new/corrected text:
#   B10 (row for BDKRB2): "9+9"           -> "+9+9"
#   B11 (row for CKM):    "985+185/117/0" -> "985+185/1170"
#   B13 (row for EPAS1 (hCV2148918)): "TT" -> "CC"
#
# B10's new value starts with "+", so Excel treats it as text entered with a
# leading apostrophe (quote-prefix) - that's what flips its cell style to the
# quotePrefix variant, centered, same font as the other B-column cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "'+9+9"
$ws.Range("B11").Value = "985+185/1170"
$ws.Range("B13").Value = "CC"

# Reflect the reviewer's new scroll position / selection on the sheet.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select()
